$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.01249966666666667
$ws.Range("H2").Value = 0.037499
$ws.Range("M2").Value = 2.423077
$ws.Range("N2").Value = 7.269231
$ws.Range("O2").Value = 0.03932631260408408
$ws.Range("P2").Value = 0.03932631260408408
$ws.Range("Q2").Value = 0.03028765480766666
$ws.Range("R2").Value = 0.272588893269
$ws.Range("S2").Value = 0.03932631260408408
$ws.Range("T2").Value = 0.03932631260408408

# Row 3
$ws.Range("G3").Value = 0.01249966666666667
$ws.Range("H3").Value = 0.037499
$ws.Range("O3").Value = 0.2611559628478186
$ws.Range("P3").Value = 0.2611559628478186
$ws.Range("Q3").Value = 0.2011325529888889
$ws.Range("R3").Value = 1.8101929769
$ws.Range("S3").Value = 0.2611559628478186
$ws.Range("T3").Value = 0.2611559628478186

# Row 4
$ws.Range("G4").Value = 0.01249966666666667
$ws.Range("H4").Value = 0.037499
$ws.Range("O4").Value = 0.6995177245480974
$ws.Range("P4").Value = 0.6995177245480974
$ws.Range("Q4").Value = 0.5387423831533333
$ws.Range("R4").Value = 4.84868144838
$ws.Range("S4").Value = 0.6995177245480974
$ws.Range("T4").Value = 0.6995177245480974
